# FEB_SN3_30kW.xlsx - "Fixed point mass derating. Two track converges in coarse grain"
#
# Updates the vehicle parameters on the "Info" sheet:
#   - Total Mass (C4):              325   -> 245
#   - Disc Outer Diameter (C15):    266.7 -> 203.2
#   - Grip Factor Multiplier (C22): 1     -> 0.9, with a new comment explaining the derating
#   - Tyre Radius (C23):            266.7 -> 203.2
# and leaves the sheet's selection on C26 (near the edited Grip Factor row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# --- Point mass derating -------------------------------------------------
$ws.Range("C4").Value = 245

# --- Two-track wheel/disc geometry converged in coarse grain -------------
$ws.Range("C15").Value = 203.2
$ws.Range("C23").Value = 203.2

# --- Grip factor multiplier derated, with justification comment ----------
$ws.Range("C22").Value = 0.9
$ws.Range("E22").Value = "Assumed 0.6-0.7 for accuracy"

# --- Leave the selection where the author left off working ---------------
$ws.Activate()
$ws.Range("C26").Select()
